$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Ilut & Valchev) entered first: B16, C16, D16, A16
$ws.Range("B16").Value = "Model of costly reasoning to update beliefs about optimal mapping of econ states to actions."
$ws.Range("C16").Value = "Key result: agents reason more about a state when its unusual --> state/history-dependent behavior with 'learning traps': endogenous familiar regions of state space where behavior appears to follow past experience based heuristics. Traps have empirically desirable properties: MPC higher, hand-to-mout more frequent and persistent and more welath inequality."
$ws.Range("D16").Value = "In an incomplete market setting, a learning-type setting produces desirable empirical features. "
$ws.Range("A16").Value = "Ilut & Valchev 2020 WP, Economic Agents as Imperfect Problem Solvers"

# Row 15 (Candia, Coibion & Goro) entered second: A15, B15, C15
$ws.Range("A15").Value = "Candia, Coibion & Goro 2020, Communication and the Beliefs of Economic Agents"
$ws.Range("B15").Value = "Look at surveys of expectations to examine how beliefs affect actions. Info provision about inflation can have opposite effects on actions depending on whether agents interpret them as supply  or demand side. As opposed to profi forecasters, HHs (but also many firms) tend to interpret inflation as supply-side, thus leading to actions that depress the economy. "
$ws.Range("C15").Value = "This may have been what Jenny mentioned. No b/c that was a JMP… and also in that JMP firms interpreted inflation as a demand shock, so it was expansionary"

# Row heights to match the other wrapped-text rows (ht=60)
$ws.Rows(15).RowHeight = 60
$ws.Rows(16).RowHeight = 60

# D16 carries the same "highlighted" style used by rows 4 and 10 (yellow fill + wrap)
$ws.Range("D16").Interior.Color = 65535
$ws.Range("D16").WrapText = $true

# Move the view / selection like the author's session: scrolled down, selection on F14
$ws.Range("F14").Select()

